# Add the new "bets" row 45 (bet #44, settled 2023-10-16) and let the
# dependent "resumen" summary formulas (which pull the latest percentage
# from bets!M:M) recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bets")

$ws.Range("C45").Value2 = 1
$ws.Range("D45").Value2 = 1683895.0160600001
$ws.Range("E45").Value2 = 100
$ws.Range("F45").Formula = "=D45+E45"
$ws.Range("G45").Value2 = "ESPORTS"
$ws.Range("H45").Value2 = "WORLDS 2023"
$ws.Range("I45").Value2 = "BLG"
$ws.Range("J45").Value2 = "ACE MAPA 1"
$ws.Range("K45").Value2 = 1
$ws.Range("L45").Value2 = 0
$ws.Range("M45").Formula = "=ROUND((F45/`$D`$31-1)*100, 3)+`$M`$29"

# Reflect where the user's cursor ended up after entering the row (one
# cell below the last edited column E, matching the saved selection).
$ws.Range("E46").Select() | Out-Null
